$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for 7f9ba23f... and c695d81b... rows
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-03 00:18:07"
$wsOverview.Range("G3").Value = "2016-09-03 00:18:07"

# "zh-cn" sheet: Priority, Correspond Handoff Datetime, Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-09-03 00:17:59"
$wsZhCn.Range("H3").Value = "2016-09-03 00:17:59"
$wsZhCn.Range("K2").Value = "2016-09-03 00:18:29"
$wsZhCn.Range("K3").Value = "2016-09-03 00:18:29"

# "de-de" sheet: Correspond Handoff Datetime, Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-03 00:18:07"
$wsDeDe.Range("H3").Value = "2016-09-03 00:18:07"
$wsDeDe.Range("K2").Value = "2016-09-03 00:18:37"
$wsDeDe.Range("K3").Value = "2016-09-03 00:18:37"
